$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest scraped
# crypto quotes. Price values are prefixed with a leading apostrophe so Excel
# keeps them as text (matching the "51.706.30"-style grouped formatting used
# by the site) instead of parsing them as numbers.

$ws.Range("D2").Value = "'51.688.50"
$ws.Range("E2").Value = '  +1.41%  '
$ws.Range("D3").Value = "'3.031.88"
$ws.Range("E3").Value = '  +2.79%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = "'380.89"
$ws.Range("E5").Value = '  +0.62%  '
$ws.Range("D6").Value = "'102.95"
$ws.Range("E6").Value = '  +1.38%  '
$ws.Range("E7").Value = '  +0.86%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = "'0.594"
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("D10").Value = "'36.83"
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("D13").Value = "'3.529.23"
$ws.Range("E13").Value = '  +3.28%  '
$ws.Range("D14").Value = "'18.58"
$ws.Range("E14").Value = '  +1.57%  '
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").Value = "'3.039.77"
$ws.Range("E16").Value = '  +3.10%  '
$ws.Range("D17").Value = "'0.978"
$ws.Range("E17").Value = '  -3.28%  '
$ws.Range("D18").Value = "'10.52"
$ws.Range("E18").Value = '  -13.70%  '
$ws.Range("D19").Value = "'51.707.81"
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("D20").Value = "'3.07"
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("E21").Value = '  +1.02%  '
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("D23").Value = "'70.11"
$ws.Range("E23").Value = '  +1.02%  '
$ws.Range("D24").Value = "'268.52"
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("D25").Value = "'3.17"
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("D26").Value = "'8.24"
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("D27").Value = "'7.58"
$ws.Range("E27").Value = '  +7.62%  '
$ws.Range("E28").Value = '  +6.39%  '
$ws.Range("E29").Value = '  +2.62%  '
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("D31").Value = "'0.109"
$ws.Range("E31").Value = '  +1.05%  '
$ws.Range("E32").Value = '  +0.96%  '
$ws.Range("D33").Value = "'34.08"
$ws.Range("E33").Value = '  +1.25%  '
$ws.Range("D34").Value = "'50.52"
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").Value = "'0.0447"
$ws.Range("E36").Value = '  +3.70%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").Value = "'3.35"
$ws.Range("E38").Value = '  +7.37%  '
$ws.Range("D39").Value = "'0.292"
$ws.Range("E39").Value = '  +13.20%  '
$ws.Range("D40").Value = "'17.11"
$ws.Range("E40").Value = '  +3.35%  '
$ws.Range("E41").Value = '  +2.84%  '
$ws.Range("E42").Value = '  +1.96%  '
$ws.Range("D43").Value = "'127.81"
$ws.Range("E43").Value = '  +8.51%  '
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("E45").Value = '  +6.47%  '
$ws.Range("D46").Value = "'21.80"
$ws.Range("E46").Value = '  +2.30%  '
$ws.Range("D47").Value = "'2.09"
$ws.Range("E47").Value = '  +3.98%  '
$ws.Range("E48").Value = '  +3.32%  '
$ws.Range("D49").Value = "'2.035.16"
$ws.Range("E49").Value = '  +1.65%  '
$ws.Range("D50").Value = "'3.335.32"
$ws.Range("E50").Value = '  +2.92%  '
$ws.Range("E51").Value = '  +0.96%  '
